# Apply commit "feat: add 2022-Q3 data"
#  1. Insert a new worksheet "2022-Q3" right after "总计" (before the existing
#     "2022-Q2" sheet) containing the fund position breakdown for the new quarter.
#  2. Insert a new summary row at the top of the "总计" data table reflecting the
#     new quarter's totals, shifting the existing quarterly summary rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q3" worksheet
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"

# Header row (values + formatting) copied from an existing quarter sheet - the
# header text/style is identical across every quarterly breakdown sheet.
$q2.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

# Column A (the 0-based row index, styled bold+centered+bordered) also matches
# the layout used on every quarter sheet, so reuse it (gives us both the style
# and the correct sequential values 0..13 in one shot).
$q2.Range("A2:A15").Copy($newSheet.Range("A2:A15"))

$q3Data = @(
    @("004666", "长城久嘉创新成长灵活配置混合A", "25.06", "84.47", "4.56", "1.1427", 7),
    @("015115", "汇添富中国高端制造股票D",       "24.58", "83.58", "2.90", "0.7128", 10),
    @("001725", "汇添富中国高端制造股票A",       "24.16", "83.58", "2.90", "0.7006", 10),
    @("010052", "长城久嘉创新成长灵活配置混合C", "6.16",  "84.47", "4.56", "0.2809", 7),
    @("161605", "融通蓝筹成长混合",               "5.16",  "73.64", "3.82", "0.1971", 4),
    @("000717", "融通转型三动力灵活配置混合A",   "3.20",  "88.91", "5.54", "0.1773", 5),
    @("560002", "益民红利成长混合",               "3.32",  "76.71", "2.58", "0.0857", 10),
    @("014606", "招商高端装备混合A",               "1.56",  "93.76", "3.41", "0.0532", 5),
    @("014607", "招商高端装备混合C",               "1.43",  "93.76", "3.41", "0.0488", 5),
    @("002681", "金鹰元和灵活配置混合A",           "0.30",  "87.46", "5.03", "0.0151", 7),
    @("002682", "金鹰元和灵活配置混合C",           "0.23",  "87.46", "5.03", "0.0116", 7),
    @("000432", "中银优秀企业混合",               "0.17",  "86.02", "2.43", "0.0041", 10),
    @("009828", "融通转型三动力灵活配置混合C",   "0.03",  "88.91", "5.54", "0.0017", 5),
    @("015114", "汇添富中国高端制造股票C",       "0.04",  "83.58", "2.90", "0.0012", 10)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $row = $i + 2
    $vals = $q3Data[$i]

    $newSheet.Cells.Item($row, 2).Value = "'" + $vals[0]
    $newSheet.Cells.Item($row, 3).Value = $vals[1]
    $newSheet.Cells.Item($row, 4).Value = "'" + $vals[2]
    $newSheet.Cells.Item($row, 5).Value = "'" + $vals[3]
    $newSheet.Cells.Item($row, 6).Value = "'" + $vals[4]
    $newSheet.Cells.Item($row, 7).Value = "'" + $vals[5]
    $newSheet.Cells.Item($row, 8).Value = $vals[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Insert a fresh row so the existing quarters shift down by one.
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

# Re-use the bold/centered/bordered style from the row below for column A.
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 14
$total.Range("D2").Value = 3.43

# The index column is purely positional (0-based row offset) - keep it in
# sequence for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
